$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for specific rows per the diff
$ws.Range("F3").Value = -1
$ws.Range("F9").Value = 4
$ws.Range("F11").Value = -4
$ws.Range("F13").Value = 3
$ws.Range("F16").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("F21").Value = -3
$ws.Range("F22").Value = 5
$ws.Range("F23").Value = -6
